$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column C ("Förändrad") holds the last-changed date for every record (rows 2-338).
# All of them advance from 2023-09-10 (serial 45179) to 2023-09-11 (serial 45180).
$ws.Range("C2:C338").Value = 45180
